# Restore cell C10 ("Rules" sheet, row for R10, "From" column) from 18 to 1,
# matching the commit's revision-restore of the Sample Project rules table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 1
